$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 6 (shifts old row 6 "hostname/adeye03u" down to row 7)
$ws.Rows.Item(6).Insert()

# Fill the new row 6 with the new point map entry
$ws.Range("B6").Value = "/opt/ros/kinetic"
$ws.Range("A6").Value = "ROS_folder"

# Select A6 as the active cell, matching the saved selection state
$ws.Range("A6").Select()
